$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$headerRange = $ws.Range("I1:J1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 4
$ws.Range("J3").Value = 5

$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 8

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
